$wb = $excel.ActiveWorkbook

# Rename shared strings used across the workbook:
#   "one.1" -> "one[1]"
#   "z.0"   -> "z[0]"
$ws1 = $wb.Worksheets.Item("s1")
$ws2 = $wb.Worksheets.Item("s2")

$ws1.Range("F1").Value = "one[1]"
$ws1.Range("G1").Value = "z[0]"

$ws2.Range("A7").Value = "one[1]"
$ws2.Range("A8").Value = "z[0]"

# Update selections / active sheet.
$ws2.Range("I7").Select()
$ws1.Select()
$ws1.Range("G11").Select()
